$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the table data (rows 2-5) with the new exposure-site rows.
$data = @(
    @("Brandon Park", "Kmart, Brandon Park Shopping Centre  Cnr Springvale Rd and Ferntree Gully Rd  Brandon Park VIC 3170", "4:35pm - 5:10pm  31/1/2021", "Case attended venue", "old"),
    @("Brighton", "North Point Cafe  2B North Rd  Brighton, VIC 3186", "8:10am - 9:30am  31/1/2021", "Case dined outside and used bathroom", "old"),
    @("Keysborough", "Kmart, Parkmore Keysborough Shopping Centre  C/317 Cheltenham Rd  Keysborough VIC 3173", "4:00pm - 5:00pm  31/1/2021", "Case attended venue", "old"),
    @("Springvale", "Coles Springvale  825 Dandenong Rd  Springvale VIC 3171", "5:00pm - 6:00pm  31/1/2021", "Case attended venue", "old")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

# Resize the columns to fit the new (wider) data, matching Excel's "best fit" autosize.
$ws.Columns.Item(1).ColumnWidth = 10.498697916666666
$ws.Columns.Item(2).ColumnWidth = 82.09635416666667
$ws.Columns.Item(3).ColumnWidth = 23.166666666666668
$ws.Columns.Item(4).ColumnWidth = 30.764322916666668

# Select all columns A:E (whole-column selection), matching the saved view state.
$ws.Range("A1:E1048576").Select() | Out-Null
